$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (rows 20..50 shift down to 21..51),
# then fill it in with a new weekly "Perejil" price observation.
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(20, 3).Value = 'Ñuble'
$ws.Cells.Item(20, 4).Value = 45044
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = 100112044
$ws.Cells.Item(20, 7).Value = 'Perejil'
$ws.Cells.Item(20, 8).Value = 'Sin especificar'
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1500
$ws.Cells.Item(20, 14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 1500
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = 'Hortaliza'
